$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: add week 37 column label, matching header style of AM1 (bold, centered)
$ws.Range("AN1").Value = "37"
$ws.Range("AN1").Font.Bold = $true
$ws.Range("AN1").HorizontalAlignment = -4108

# New week-37 case counts + a few corrected historical values
$ws.Range("AN2").Value = 0
$ws.Range("AN3").Value = 0
$ws.Range("AN5").Value = 0
$ws.Range("AN6").Value = 33
$ws.Range("AN7").Value = 2
$ws.Range("AN8").Value = 16
$ws.Range("AN9").Value = 0
$ws.Range("AN10").Value = 0
$ws.Range("AN12").Value = 0
$ws.Range("AN14").Value = 0
$ws.Range("AN15").Value = 0
$ws.Range("AN16").Value = 0
$ws.Range("AN17").Value = 0
$ws.Range("AN21").Value = 0
$ws.Range("AN23").Value = 0
$ws.Range("AN25").Value = 2
$ws.Range("AN26").Value = 0
$ws.Range("AA28").Value = 1
$ws.Range("AB28").Value = 5
$ws.Range("AD28").Value = 3
$ws.Range("AE28").Value = 0
$ws.Range("AF28").Value = 3
$ws.Range("AG28").Value = 2
$ws.Range("AH28").Value = 0
$ws.Range("AI28").Value = 1
$ws.Range("AJ28").Value = 0
$ws.Range("AK28").Value = 1
$ws.Range("AL28").Value = 3
$ws.Range("AN28").Value = 3
$ws.Range("T28").Value = 1
$ws.Range("U28").Value = 1
$ws.Range("V28").Value = 1
$ws.Range("X28").Value = 2
$ws.Range("Z28").Value = 2
$ws.Range("AN29").Value = 2
$ws.Range("AN30").Value = 13
$ws.Range("AM31").Value = 0
$ws.Range("AN31").Value = 0
$ws.Range("AM35").Value = 13
$ws.Range("AN35").Value = 11
$ws.Range("M35").Value = 13
$ws.Range("AN36").Value = 0
$ws.Range("AN37").Value = 0
$ws.Range("AN38").Value = 0
$ws.Range("AN41").Value = 0
$ws.Range("AN42").Value = 0
$ws.Range("AM43").Value = 0
$ws.Range("AN43").Value = 0
$ws.Range("AN44").Value = 0
$ws.Range("AN45").Value = 0
$ws.Range("AN46").Value = 0
$ws.Range("AN47").Value = 0
$ws.Range("AN48").Value = 0
$ws.Range("AN49").Value = 0
$ws.Range("AN50").Value = 0
$ws.Range("AN51").Value = 0
$ws.Range("AM52").Value = 0
$ws.Range("AN53").Value = 0
$ws.Range("AN54").Value = 0
$ws.Range("AN55").Value = 0
$ws.Range("AN56").Value = 0
$ws.Range("AN57").Value = 0
$ws.Range("AN58").Value = 0
